$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Student Profile")

# Update student's first / last name
$ws.Range("B3").Value = "Abeer"
$ws.Range("B4").Value = "Bajpai"

# Accomplishments: switch checked option from "Student on the year" (B20)
# to "One of the Top Students in my class" (B21)
$ws.Range("B20").Value = ""
$ws.Range("B21").Value = "X"

# Positive Personality Traits: switch checked option from "enthusiastic" (B36)
# to "bright" (B32)
$ws.Range("B32").Value = "X"
$ws.Range("B36").Value = ""

# Academic Skills: switch checked option from "great presentation skills" (B78)
# to "problem-solving skills" (B70)
$ws.Range("B70").Value = "X"
$ws.Range("B78").Value = ""

# Restore the active cell selection on the sheet
$ws.Range("B9").Select()
